$wb = $excel.ActiveWorkbook

# --- Sheet2 "RM": zero-out the reserve margin values in row 2 (B2:AK2) ---
$wsRM = $wb.Worksheets.Item("RM")
for ($col = 2; $col -le 37; $col++) {
    $wsRM.Cells.Item(2, $col).Value = 0
}

# --- Sheet1 "About": add the new EPS 3.0.0 note rows ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A14").Value = "U.S. EPS 3.0.0 Note"
$wsAbout.Range("A14").Style = "Normal"
$wsAbout.Range("A14").Font.Bold = $true

$wsAbout.Range("A15").Value = "In EPS 3.0.0, a new peaking calculation approach was adopted.  The new approach uses Equipment Load Factors (ELFs)"
$wsAbout.Range("A16").Value = "that may already account for some or all of the reserve margin that utilities consider.  We set the reservere margin to"
$wsAbout.Range("A17").Value = "zero, but we leave the data reference above in place in case we wish to use a non-zero reserve margin in the future."
